# "Esthetic modifications of the template"
#  - tighten page margins (top/right/bottom/left) from 0.5" (720 twips)
#    to 36.85pt (737 twips); header/footer/gutter stay as-is
#  - make the built-in "Normal" style justify its paragraphs (jc=both)

$d = $word.ActiveDocument

# --- page margins -----------------------------------------------------
$ps = $d.Sections(1).PageSetup
$ps.TopMargin    = 36.85
$ps.BottomMargin = 36.85
$ps.LeftMargin   = 36.85
$ps.RightMargin  = 36.85

# --- Normal style paragraph alignment ---------------------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.Alignment = 3   # wdAlignParagraphJustify
